$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1674.6666
$ws.Range("J17").Value = 1674.6666
$ws.Range("L17").Value = 5023.9998
$ws.Range("N17").Value = -5359.9998

$ws.Range("H53").Value = 13334156
$ws.Range("J53").Value = 998.46155
$ws.Range("L53").Value = 998.46155
$ws.Range("N53").Value = -2272.46155

$ws.Range("H62").Value = 15628856
$ws.Range("I62").Value = 17860406
$ws.Range("K62").Value = 17860406
$ws.Range("M62").Value = -17859782

$ws.Range("H65").Value = 15628856
$ws.Range("I65").Value = 17860406
$ws.Range("K65").Value = 89302030
$ws.Range("M65").Value = -89298910

$ws.Range("H69").Value = 17500
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()

$ws.Range("H72").Value = 17500
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()

$ws.Range("H76").Value = 125006100
$ws.Range("I76").Value = 6599
$ws.Range("J76").Value = 250005600
$ws.Range("K76").Value = 6599
$ws.Range("L76").Value = 250005600
$ws.Range("M76").Value = -6284
$ws.Range("N76").Value = -250006230

$ws.Range("H79").Value = 125006100
$ws.Range("I79").Value = 6599
$ws.Range("J79").Value = 250005600
$ws.Range("K79").Value = 6599
$ws.Range("L79").Value = 250005600
$ws.Range("M79").Value = -5507
$ws.Range("N79").Value = -250007784

$ws.Range("H86").Value = 3762175.8
$ws.Range("J86").Value = 7520925
$ws.Range("L86").Value = 7520925
$ws.Range("N86").Value = -7523171

$ws.Range("H89").Value = 3762175.8
$ws.Range("J89").Value = 7520925
$ws.Range("L89").Value = 37604625
$ws.Range("N89").Value = -37615857

$ws.Range("H100").Value = 9348.941000000001
$ws.Range("I100").Value = 2156.1428
$ws.Range("J100").Value = 14383.9
$ws.Range("K100").Value = 2156.1428
$ws.Range("L100").Value = 14383.9
$ws.Range("M100").Value = -1615.1428
$ws.Range("N100").Value = -15465.9

$ws.Range("H111").Value = 95361.27
$ws.Range("I111").Value = 128623.125
$ws.Range("J111").Value = 6663
$ws.Range("K111").Value = 385869.375
$ws.Range("L111").Value = 19989
$ws.Range("M111").Value = -382802.375
$ws.Range("N111").Value = -26123

$ws.Range("H132").Value = 2621.5615
$ws.Range("I132").Value = 2243.258
$ws.Range("J132").Value = 4753.8184
$ws.Range("K132").Value = 6729.773999999999
$ws.Range("L132").Value = 14261.4552
$ws.Range("M132").Value = -4199.773999999999
$ws.Range("N132").Value = -19321.4552

$ws.Range("H137").Value = 411002.97
$ws.Range("I137").Value = 251991.05
$ws.Range("J137").Value = 1117722.6
$ws.Range("K137").Value = 755973.1499999999
$ws.Range("L137").Value = 3353167.8
$ws.Range("M137").Value = -753423.1499999999
$ws.Range("N137").Value = -3358267.8

$ws.Range("H138").Value = 4015.9866
$ws.Range("J138").Value = 4800.3726
$ws.Range("L138").Value = 14401.1178
$ws.Range("N138").Value = -24681.1178

$ws.Range("H141").Value = 1707.039
$ws.Range("I141").Value = 780.51514
$ws.Range("J141").Value = 7266.1816
$ws.Range("K141").Value = 2341.54542
$ws.Range("L141").Value = 21798.5448
$ws.Range("M141").Value = 2838.45458
$ws.Range("N141").Value = -32158.5448

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5066.64
$ws.Range("I32").Value = 4596.289
$ws.Range("K32").Value = 4596.289
$ws.Range("M32").Value = -4309.289

$ws.Range("H61").Value = 1333.54
$ws.Range("I61").Value = 753.84784
$ws.Range("K61").Value = 753.84784
$ws.Range("M61").Value = -541.84784

$ws.Range("H74").Value = 3429.111
$ws.Range("I74").Value = 3378.3333
$ws.Range("J74").Value = 3530.6667
$ws.Range("K74").Value = 3378.3333
$ws.Range("L74").Value = 3530.6667
$ws.Range("M74").Value = -2504.3333
$ws.Range("N74").Value = -5278.6667

$ws.Range("H77").Value = 3429.111
$ws.Range("I77").Value = 3378.3333
$ws.Range("J77").Value = 3530.6667
$ws.Range("K77").Value = 16891.6665
$ws.Range("L77").Value = 17653.3335
$ws.Range("M77").Value = -12523.6665
$ws.Range("N77").Value = -26389.3335

$ws.Range("H136").Value = 1333.54
$ws.Range("I136").Value = 753.84784
$ws.Range("K136").Value = 2261.54352
$ws.Range("M136").Value = 288.4564799999998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 17936.797
$ws.Range("I134").Value = 2389.44
$ws.Range("K134").Value = 7168.32
$ws.Range("M134").Value = -4633.32

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 470006.34
$ws.Range("I31").Value = 272831.84
$ws.Range("K31").Value = 272831.84
$ws.Range("M31").Value = -272536.84

$ws.Range("H34").Value = 470006.34
$ws.Range("I34").Value = 272831.84
$ws.Range("K34").Value = 272831.84
$ws.Range("M34").Value = -272629.84

$ws.Range("H59").Value = 31705
$ws.Range("J59").Value = 31705
$ws.Range("L59").Value = 31705
$ws.Range("N59").Value = -33995

$ws.Range("H68").Value = 120000
$ws.Range("J68").Value = 120000
$ws.Range("L68").Value = 120000
$ws.Range("N68").Value = -121498

$ws.Range("H71").Value = 120000
$ws.Range("J71").Value = 120000
$ws.Range("L71").Value = 360000
$ws.Range("N71").Value = -367488

$ws.Range("H74").Value = 92466.664
$ws.Range("J74").Value = 92466.664
$ws.Range("L74").Value = 92466.664
$ws.Range("N74").Value = -94214.664

$ws.Range("H77").Value = 92466.664
$ws.Range("J77").Value = 92466.664
$ws.Range("L77").Value = 277399.992
$ws.Range("N77").Value = -286135.992

$ws.Range("H107").Value = 797.4375
$ws.Range("I107").Value = 724.1905
$ws.Range("J107").Value = 937.2727
$ws.Range("K107").Value = 724.1905
$ws.Range("L107").Value = 937.2727
$ws.Range("M107").Value = 1195.8095
$ws.Range("N107").Value = -4777.2727

$ws.Range("H132").Value = 2372.6191
$ws.Range("I132").Value = 1749.2898
$ws.Range("K132").Value = 5247.8694
$ws.Range("M132").Value = -2717.8694

$ws.Range("H134").Value = 208370.47
$ws.Range("I134").Value = 117999.74
$ws.Range("K134").Value = 353999.22
$ws.Range("M134").Value = -351464.22

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 6957.737
$ws.Range("I56").Value = 6957.737
$ws.Range("K56").Value = 6957.737
$ws.Range("M56").Value = -6427.737

$ws.Range("H134").Value = 1051.8
$ws.Range("I134").Value = 1051.8
$ws.Range("K134").Value = 3155.4
$ws.Range("M134").Value = 1914.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 18287144
$ws.Range("I11").Value = 5752501
$ws.Range("K11").Value = 5752501
$ws.Range("M11").Value = -5752362

$ws.Range("H12").Value = 1000
$ws.Range("I12").Value = 1000
$ws.Range("K12").Value = 1000
$ws.Range("M12").Value = -860

$ws.Range("I14").Value = 58021630
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 58021630
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -58021462
$ws.Range("N14").ClearContents()

$ws.Range("H132").Value = 177009.08
$ws.Range("I132").Value = 201755.55
$ws.Range("K132").Value = 605266.6499999999
$ws.Range("M132").Value = -602736.6499999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1995.4
$ws.Range("I16").Value = 1997.6666
$ws.Range("K16").Value = 1997.6666
$ws.Range("M16").Value = -1827.6666

$ws.Range("H17").Value = 3009
$ws.Range("J17").Value = 3009
$ws.Range("L17").Value = 3009
$ws.Range("N17").Value = -3349

$ws.Range("H18").Value = 6669
$ws.Range("I18").Value = 7500.5
$ws.Range("J18").Value = 5006
$ws.Range("K18").Value = 7500.5
$ws.Range("L18").Value = 5006
$ws.Range("M18").Value = -7328.5
$ws.Range("N18").Value = -5350

$ws.Range("H93").Value = 3051.9285
$ws.Range("I93").Value = 2757.8
$ws.Range("J93").Value = 3787.25
$ws.Range("K93").Value = 2757.8
$ws.Range("L93").Value = 3787.25
$ws.Range("M93").Value = -1509.8
$ws.Range("N93").Value = -6283.25

$ws.Range("H100").Value = 3447.5
$ws.Range("I100").Value = 2515.2
$ws.Range("K100").Value = 2515.2
$ws.Range("M100").Value = -1974.2

$ws.Range("H132").Value = 3809.3208
$ws.Range("I132").Value = 3287.6
$ws.Range("J132").Value = 5414.615
$ws.Range("K132").Value = 9862.799999999999
$ws.Range("L132").Value = 16243.845
$ws.Range("M132").Value = -7332.799999999999
$ws.Range("N132").Value = -21303.845

$ws.Range("H136").Value = 395010.5
$ws.Range("I136").Value = 478155.56
$ws.Range("J136").Value = 7000.3335
$ws.Range("K136").Value = 1434466.68
$ws.Range("L136").Value = 21001.0005
$ws.Range("M136").Value = -1431916.68
$ws.Range("N136").Value = -26101.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1666.6666
$ws.Range("I81").Value = 1666.6666
$ws.Range("K81").Value = 3333.3332
$ws.Range("M81").Value = -2272.3332

$ws.Range("H84").Value = 1666.6666
$ws.Range("I84").Value = 1666.6666
$ws.Range("K84").Value = 16666.666
$ws.Range("M84").Value = -11362.666

$ws.Range("H107").Value = 1013
$ws.Range("I107").Value = 919.75
$ws.Range("J107").Value = 1199.5
$ws.Range("K107").Value = 2759.25
$ws.Range("L107").Value = 3598.5
$ws.Range("M107").Value = -839.25
$ws.Range("N107").Value = -7438.5

$ws.Range("H113").Value = 854.1905
$ws.Range("I113").Value = 576.5333000000001
$ws.Range("J113").Value = 1548.3334
$ws.Range("K113").Value = 1729.5999
$ws.Range("L113").Value = 4645.0002
$ws.Range("M113").Value = 440.4000999999998
$ws.Range("N113").Value = -8985.0002

$ws.Range("H122").Value = 25644686
$ws.Range("I122").Value = 55557344
$ws.Range("J122").Value = 5267.7144
$ws.Range("K122").Value = 166672032
$ws.Range("L122").Value = 15803.1432
$ws.Range("M122").Value = -166669582
$ws.Range("N122").Value = -20703.1432
